$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("perf_stats")

# Row 10 - Skew (log returns)
$ws.Range("B10").Value = -0.78
$ws.Range("C10").Value = -0.64
$ws.Range("D10").Value = -0.59
$ws.Range("E10").Value = -1.51
$ws.Range("F10").Value = -1.59
$ws.Range("G10").Value = 0.06
$ws.Range("H10").Value = -0.63
$ws.Range("I10").Value = -1.36

# Row 11 - Kurtosis (log returns)
$ws.Range("B11").Value = 4.64
$ws.Range("C11").Value = 3.88
$ws.Range("D11").Value = 3.68
$ws.Range("E11").Value = 14.07
$ws.Range("F11").Value = 14.4
$ws.Range("G11").Value = 19.27
$ws.Range("H11").Value = 4.1
$ws.Range("I11").Value = 11.7
